# Fruta / hortaliza, semanal
# The D (Fecha), J (Volumen), K (Precio mínimo), L (Precio máximo),
# M (Precio promedio ponderado) and P (Precio $/Kg) values get shuffled
# across the data rows (2-16); every other column is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2..16, in order: D, J, K, L, M, P
$newData = @{
  2  = @(44547, 200, 13000, 14000, 13500, 750)
  3  = @(45142, 400, 17000, 18000, 17500, 972)
  4  = @(44977, 400, 16500, 17000, 16750, 931)
  5  = @(45068, 400, 16000, 17000, 16500, 917)
  6  = @(45154, 500, 16500, 17000, 16750, 931)
  7  = @(44568, 500, 15000, 16000, 15500, 861)
  8  = @(44984, 200, 17000, 18000, 17500, 972)
  9  = @(44998, 320, 17000, 18000, 17500, 972)
  10 = @(45005, 200, 17000, 18000, 17500, 972)
  11 = @(45152, 500, 16000, 17000, 16500, 917)
  12 = @(44957, 400, 21000, 22000, 21500, 1194)
  13 = @(44960, 400, 19500, 20000, 19750, 1097)
  14 = @(45117, 300, 17000, 18000, 17500, 972)
  15 = @(44557, 400, 13000, 14000, 13500, 750)
  16 = @(44964, 300, 20000, 21000, 20500, 1139)
}

foreach ($row in $newData.Keys) {
  $vals = $newData[$row]
  $ws.Cells.Item($row, 4).Value  = $vals[0]  # D - Fecha
  $ws.Cells.Item($row, 10).Value = $vals[1]  # J - Volumen
  $ws.Cells.Item($row, 11).Value = $vals[2]  # K - Precio minimo
  $ws.Cells.Item($row, 12).Value = $vals[3]  # L - Precio maximo
  $ws.Cells.Item($row, 13).Value = $vals[4]  # M - Precio promedio ponderado
  $ws.Cells.Item($row, 16).Value = $vals[5]  # P - Precio $/Kg
}
